$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.101656556129456
$ws.Range("B1").Value = 4.115572929382324
$ws.Range("C1").Value = 6.046109199523926
$ws.Range("D1").Value = 1.486738324165344
$ws.Range("E1").Value = 0.8377549052238464
